$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.2033639143730887
$ws.Range("C2").Value = 0.5565749235474006
$ws.Range("J2").Value = 0.01070336391437309
$ws.Range("P2").Value = 0.1498470948012232
$ws.Range("S2").Value = 0.07951070336391437
$ws.Range("B3").Value = 0.01025641025641026
$ws.Range("C3").Value = 0.02307692307692308
$ws.Range("J3").Value = 0.05128205128205128
$ws.Range("P3").Value = 0.6820512820512821
$ws.Range("S3").Value = 0.2333333333333333
$ws.Range("J4").Value = 0.08163265306122448
$ws.Range("P4").Value = 0.6020408163265306
$ws.Range("S4").Value = 0.3163265306122449
$ws.Range("B6").Value = 0.05580357142857143
$ws.Range("D6").Value = 0.008928571428571428
$ws.Range("F6").Value = 0.05803571428571429
$ws.Range("J6").Value = 0.2700892857142857
$ws.Range("O6").Value = 0.03125
$ws.Range("Q6").Value = 0.1316964285714286
$ws.Range("R6").Value = 0.07589285714285714
$ws.Range("S6").Value = 0.3683035714285715
$ws.Range("B7").Value = 0.1206434316353887
$ws.Range("D7").Value = 0.005361930294906166
$ws.Range("F7").Value = 0.05630026809651475
$ws.Range("J7").Value = 0.1689008042895442
$ws.Range("O7").Value = 0.01876675603217158
$ws.Range("Q7").Value = 0.1554959785522788
$ws.Range("R7").Value = 0.09651474530831099
$ws.Range("S7").Value = 0.3780160857908847
$ws.Range("B8").Value = 0.1135073779795687
$ws.Range("D8").Value = 0.02610669693530079
$ws.Range("E8").Value = 0.001135073779795687
$ws.Range("F8").Value = 0.06242905788876277
$ws.Range("J8").Value = 0.1112372304199773
$ws.Range("O8").Value = 0.02043132803632236
$ws.Range("Q8").Value = 0.1736662883087401
$ws.Range("R8").Value = 0.09534619750283768
$ws.Range("S8").Value = 0.3961407491486947
$ws.Range("B9").Value = 0.1014925373134328
$ws.Range("D9").Value = 0.0208955223880597
$ws.Range("F9").Value = 0.06567164179104477
$ws.Range("J9").Value = 0.1373134328358209
$ws.Range("O9").Value = 0.02985074626865672
$ws.Range("Q9").Value = 0.182089552238806
$ws.Range("R9").Value = 0.07761194029850746
$ws.Range("S9").Value = 0.3850746268656716
$ws.Range("B10").Value = 0.1209265802905379
$ws.Range("D10").Value = 0.02552021986650962
$ws.Range("E10").Value = 0.002748331370239498
$ws.Range("F10").Value = 0.07302709069493522
$ws.Range("J10").Value = 0.1315272870043188
$ws.Range("O10").Value = 0.01845308205732234
$ws.Range("Q10").Value = 0.1923831959167648
$ws.Range("R10").Value = 0.07930899096976836
$ws.Range("S10").Value = 0.3561052218296035
$ws.Range("F11").Value = 0.001821493624772313
$ws.Range("G11").Value = 0.1475409836065574
$ws.Range("J11").Value = 0.1038251366120219
$ws.Range("K11").Value = 0.1912568306010929
$ws.Range("L11").Value = 0.5409836065573771
$ws.Range("S11").Value = 0.01457194899817851
$ws.Range("G12").Value = 0.7682539682539683
$ws.Range("J12").Value = 0.1492063492063492
$ws.Range("K12").Value = 0.02222222222222222
$ws.Range("L12").Value = 0.04126984126984127
$ws.Range("S12").Value = 0.01904761904761905
$ws.Range("G13").Value = 0.6555555555555556
$ws.Range("J13").Value = 0.3
$ws.Range("S13").Value = 0.04444444444444445
$ws.Range("F15").Value = 0.02262443438914027
$ws.Range("H15").Value = 0.1742081447963801
$ws.Range("I15").Value = 0.07239819004524888
$ws.Range("J15").Value = 0.3529411764705883
$ws.Range("K15").Value = 0.07466063348416289
$ws.Range("M15").Value = 0.002262443438914027
$ws.Range("N15").Value = 0.004524886877828055
$ws.Range("O15").Value = 0.04524886877828054
$ws.Range("S15").Value = 0.251131221719457
$ws.Range("F16").Value = 0.01686746987951807
$ws.Range("H16").Value = 0.1662650602409639
$ws.Range("I16").Value = 0.06265060240963856
$ws.Range("J16").Value = 0.4240963855421687
$ws.Range("K16").Value = 0.0891566265060241
$ws.Range("M16").Value = 0.03373493975903614
$ws.Range("O16").Value = 0.06746987951807229
$ws.Range("S16").Value = 0.1397590361445783
$ws.Range("F17").Value = 0.0171990171990172
$ws.Range("H17").Value = 0.2137592137592138
$ws.Range("I17").Value = 0.08476658476658476
$ws.Range("J17").Value = 0.4066339066339066
$ws.Range("K17").Value = 0.1056511056511057
$ws.Range("M17").Value = 0.01842751842751843
$ws.Range("O17").Value = 0.06511056511056511
$ws.Range("S17").Value = 0.08845208845208845
$ws.Range("F18").Value = 0.0291005291005291
$ws.Range("H18").Value = 0.1693121693121693
$ws.Range("I18").Value = 0.0873015873015873
$ws.Range("J18").Value = 0.4232804232804233
$ws.Range("K18").Value = 0.0873015873015873
$ws.Range("M18").Value = 0.02116402116402116
$ws.Range("O18").Value = 0.08201058201058201
$ws.Range("S18").Value = 0.1005291005291005
$ws.Range("F19").Value = 0.02001668056713928
$ws.Range("H19").Value = 0.2093411175979983
$ws.Range("I19").Value = 0.07422852376980818
$ws.Range("J19").Value = 0.3882402001668057
$ws.Range("K19").Value = 0.103419516263553
$ws.Range("M19").Value = 0.02460383653044203
$ws.Range("N19").Value = 0.0004170141784820684
$ws.Range("O19").Value = 0.06797331109257715
$ws.Range("S19").Value = 0.1117597998331943
